$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[4.353287029620349, 8.720442679070402]"
$ws.Range("M2").Value = 0.000000009504519393743749
$ws.Range("N2").Value = 0.0000000190090387874875
$ws.Range("P2").Value = "[-1.9497371824080805, -1.1195265111891555]"
$ws.Range("Q2").Value = 0.000000000002543520949416234
$ws.Range("R2").Value = 0.000000000002543520949416234
$ws.Range("T2").Value = "[7.547360996526498, 10.417519345055753]"
$ws.Range("X2").Value = 4.572052052052153
$ws.Range("Y2").Value = 7.962562562562742

# Row 3 updates
$ws.Range("L3").Value = "[4.272237551227706, 9.582012377639872]"
$ws.Range("M3").Value = 0.0000005216981566746881
$ws.Range("N3").Value = 0.0000005216981566746881
$ws.Range("T3").Value = "[7.422025299811736, 10.51150928351736]"
